$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the excess "1" values in column I for rows 7 and 14-18
#    (the DTR summary used to over-count a column that should be blank).
# ---------------------------------------------------------------------------
$ws.Range("I7").ClearContents()
$ws.Range("I14").ClearContents()
$ws.Range("I15").ClearContents()
$ws.Range("I16").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("I18").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add a "Legends:" header above the legend table (row 24, columns E:P)
#    using the same bold/underlined/size-15 look as the document title.
# ---------------------------------------------------------------------------
$legendHeader = $ws.Range("E24:P24")
$legendHeader.Value = " "
$ws.Range("E24").Value = "Legends:"
$legendHeader.Font.Name = "Arial"
$legendHeader.Font.Size = 15
$legendHeader.Font.Bold = $true
$legendHeader.Font.Underline = $true
$legendHeader.Merge()

# ---------------------------------------------------------------------------
# 3. Build the 3-row legend (color swatch in column E, description in F:P)
#    Row 25-26: blue swatch / late-or-undertime remark
#    Row 27-28: orange swatch / half-day remark
#    Row 29-30: red swatch / absent remark
# ---------------------------------------------------------------------------

function Set-LegendRow($swatchRange, $swatchColor, $textCell, $textRange, $text) {
    $swatchRange.Value = " "
    $swatchRange.Interior.Color = $swatchColor
    $textRange.Value = " "
    $textCell.Value = $text
    $textRange.Font.Name = "Arial"
    $textRange.Font.Size = 11
    $textRange.Font.Bold = $true
    $textRange.Font.Underline = $true
}

# Blue legend (rows 25:26)
Set-LegendRow $ws.Range("E25:E26") 13411113 $ws.Range("F25") $ws.Range("F25:P26") "Employee has request(s)/remark(s) for that day.`n*May incur late and/or undertime depending on his or her time-in and time-out."

# Orange legend (rows 27:28)
Set-LegendRow $ws.Range("E27:E28") 6737151 $ws.Range("F27") $ws.Range("F27:P28") "Employee is considered half-day because of his time-in or time-out."

# Red legend (rows 29:30)
Set-LegendRow $ws.Range("E29:E30") 6184671 $ws.Range("F29") $ws.Range("F29:P30") "Employee has no time-in and therefore, considered as absent."

# Merge the swatch / description areas now that values & formatting are set.
$ws.Range("E25:E26").Merge()
$ws.Range("F25:P26").Merge()
$ws.Range("E27:E28").Merge()
$ws.Range("F27:P28").Merge()
$ws.Range("E29:E30").Merge()
$ws.Range("F29:P30").Merge()
